$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Robert Jones row (row 10):
# Project Number: 0000/0011 -> 9002/0011
$ws.Range("B10").Value = "9002/0011"
# Member Country: India -> France
$ws.Range("P10").Value = "France"

# Update Member Country for DeLuca (row 2): USA -> BELGIUM
$ws.Range("P2").Value = "BELGIUM"

# Update the sheet view: scroll to show column L, and select P3
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("P3").Select()
